$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.017.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.599.29'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.32%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("E6").Value = '  -2.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3781'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.89%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3643'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.28%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.74'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.85%  '

$ws.Range("E10").Value = '  -4.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.002'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("E12").Value = '  -2.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.51'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.590'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.347'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.35%  '

$ws.Range("E16").Value = '  -4.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.605.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06826'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.45%  '

$ws.Range("E20").Value = '  -5.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.540'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.58%  '

$ws.Range("B22").Value = 'BitDAO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.5579'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.33%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.02%  '

$ws.Range("B25").Value = 'WrappedBTC'
$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '23.011.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.42%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.360'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.72%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.813'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.79%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.24%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '150.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.66%  '

$ws.Range("B30").Value = 'HuobiToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.230'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.34%  '

$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.62%  '

$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.329'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.00%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.818'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -12.12%  '

$ws.Range("B34").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C34").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.785.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.98%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9596'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.87%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.07562'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.64%  '

$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.31'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.70%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.248'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.39%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02701'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.11%  '

$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2527'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.37%  '

$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.08873'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.98%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.364'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.51%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7027'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.67%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.38'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.17%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.77%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6616'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.37%  '

$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.292'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.24%  '

$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.991'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.77%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.95%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07904'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.72%  '
